# Correction in SA algorithm and 746 logs: update the Fitness (column C)
# values in run_8.xlsx's log sheet to reflect the corrected run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-7  (Generation 0-5)   -> 7704
$ws.Range("C2:C7").Value = 7704

# Rows 8-84 (Generation 6-82)  -> 7310
$ws.Range("C8:C84").Value = 7310

# Rows 85-123 (Generation 83-121) -> 7293
$ws.Range("C85:C123").Value = 7293
